{"js": "// Split the paragraph containing \"{m:userdoc 'zone1'}\" (currently stored as\n// two runs: \"{m\" and \":userdoc 'zone1'}\") into four separate runs:\n// \"{\", \"m\", \":userdoc 'zone1'\", \"}\".\n// This mirrors the TokenIteratorFieldRewriterSplit parser change: field\n// delimiters (\"{\", \"}\") and the field-type token (\"m\") are now emitted as\n// their own runs instead of being glued to neighbouring text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the target paragraph by its full text content.\nconst targetText = \"{m:userdoc 'zone1'}\";\nlet target = null;\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text === targetText) {\n    target = paragraph;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find paragraph with text: \" + targetText);\n}\n\n// Grab the paragraph's own OOXML so we can preserve its attributes\n// (e.g. w:rsidP/w:rsidR/w:rsidRDefault) on the rewritten paragraph.\nconst fullRange = target.getRange();\nconst ooxmlResult = fullRange.getOoxml();\nawait context.sync();\n\nconst sourceOoxml = ooxmlResult.value;\nconst openTagMatch = sourceOoxml.match(/<w:p\\b([^>]*)>/);\nlet paragraphAttrs = openTagMatch ? openTagMatch[1] : \"\";\n// Drop volatile, auto-generated identifiers that are not part of the\n// intended edit (Word mints new ones on every OOXML round-trip anyway).\nparagraphAttrs = paragraphAttrs\n  .replace(/\\s*w14:paraId=\"[^\"]*\"/, \"\")\n  .replace(/\\s*w14:textId=\"[^\"]*\"/, \"\");\n\nconst replacementOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p' + paragraphAttrs + '>' +\n  '<w:r><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>m</w:t></w:r>' +\n  '<w:r><w:t>:userdoc \\'zone1\\'</w:t></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nfullRange.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Split the paragraph containing \"{m:userdoc 'zone1'}\" (currently stored as\n# two runs: \"{m\" and \":userdoc 'zone1'}\") into four separate runs:\n# \"{\", \"m\", \":userdoc 'zone1'\", \"}\".\n# This mirrors the TokenIteratorFieldRewriterSplit parser change: field\n# delimiters (\"{\", \"}\") and the field-type token (\"m\") are now emitted as\n# their own runs instead of being glued to neighbouring text.\n\n$d = $word.ActiveDocument\n\n$targetText = \"{m:userdoc 'zone1'}\"\n\n$target = $null\nforeach ($para in $d.Paragraphs) {\n  # Paragraph.Range.Text includes the trailing paragraph-mark character(s)\n  # (CR, or CR+BEL for list items) - trim them before comparing.\n  $t = $para.Range.Text.TrimEnd([char]13, [char]7)\n  if ($t -eq $targetText) {\n    $target = $para\n    break\n  }\n}\n\nif ($null -eq $target) {\n  throw \"Could not find paragraph with text: $targetText\"\n}\n\n$rng = $target.Range\n\n# Preserve the paragraph's own attributes (e.g. w:rsidP/w:rsidR/w:rsidRDefault)\n# on the rewritten paragraph by reading its current OOXML first.\n$srcXml = $rng.WordOpenXML\nif ($srcXml -match '<w:p\\b([^>]*)>') {\n  $pAttrs = $Matches[1]\n} else {\n  $pAttrs = \"\"\n}\n# Drop volatile, auto-generated identifiers that are not part of the\n# intended edit (Word mints new ones on every OOXML round-trip anyway).\n$pAttrs = $pAttrs -replace ' w14:paraId=\"[^\"]*\"', '' -replace ' w14:textId=\"[^\"]*\"', ''\n\n$newXml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n<pkg:xmlData>\n<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n<w:body>\n<w:p$pAttrs>\n<w:r><w:t>{</w:t></w:r>\n<w:r><w:t>m</w:t></w:r>\n<w:r><w:t>:userdoc 'zone1'</w:t></w:r>\n<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>\n</w:p>\n</w:body>\n</w:document>\n</pkg:xmlData>\n</pkg:part>\n</pkg:package>\n\"@\n\n[void]$rng.InsertXML($newXml)\n"}
